# Apply the "Added notes to tracker file" edit to the Trial Tracker workbook.
#
# Summary of the change:
#   - Trial 5 (row 9) network architecture note reverted from the "10/10"
#     description back to the "8/8" description (matches rows 5-8).
#   - The "Notes & Observations" column (Q) gains a Text number format and
#     word-wrap for the header + every data row, and rows 5-9 get a short
#     free-text note describing what was learned for that trial run.
#   - Row heights for rows 5-9 grow to fit the new wrapped text.
#   - Final selection lands on Q9 (the last cell touched by the edit).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$nl = [char]10

# --- Trial 5 (row 9): put the architecture note back to 8/8 (same text as
#     rows 5-8) instead of the 10/10 text it had before.
$ws.Range("C9").Value = "* 3-Layer NN" + $nl + "* n1 = 8" + $nl + "* n2 = 8"

# --- Notes & Observations column formatting: Text number format + wrap,
#     applied cell-by-cell (wrap first, then number format) so each cell
#     keeps its own existing font/alignment instead of picking up a new one.
$ws.Range("Q3").WrapText = $true
$ws.Range("Q3").NumberFormat = "@"

$ws.Range("Q4").WrapText = $true
$ws.Range("Q4").NumberFormat = "@"

$ws.Range("Q5").WrapText = $true
$ws.Range("Q5").NumberFormat = "@"

$ws.Range("Q6").WrapText = $true
$ws.Range("Q6").NumberFormat = "@"

$ws.Range("Q7").WrapText = $true
$ws.Range("Q7").NumberFormat = "@"

$ws.Range("Q8").WrapText = $true
$ws.Range("Q8").NumberFormat = "@"

$ws.Range("Q9").WrapText = $true
$ws.Range("Q9").NumberFormat = "@"

# --- New free-text notes for each trial.
$ws.Range("Q5").Value = "Increased the second hidden layer to 8 nodes. Overall accuracy went up in training and testing"
$ws.Range("Q6").Value = "Decreased the learning rate by an order of magnitude and it had terrible consequences"
$ws.Range("Q7").Value = "Doubled the batch size and returned the learning rate back. The increase in batch size had no good affect"
$ws.Range("Q8").Value = "Set validation_split=0.25 and got a good bumb in training and test accuracy"
$ws.Range("Q9").Value = "Increased the number of nodes in the two hidden layers to 10 and it got me over the line"

# --- Row heights grow to fit the wrapped note text.
$ws.Rows("5").RowHeight = 63.75
$ws.Rows("6").RowHeight = 57
$ws.Rows("7").RowHeight = 71.25
$ws.Rows("8").RowHeight = 52.5
$ws.Rows("9").RowHeight = 57

# --- Final selection, matching where the author's cursor ended up.
$ws.Range("Q9").Select() | Out-Null
